# "added fill form and save to folder"
# A new daily form submission happened on 2024-07-12: the "current" summary
# sheet gets a new totals row, and a brand-new daily sheet "2024-07-12" is
# appended (copy of the per-visit form layout used by the other daily sheets).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "current" summary sheet: append the totals row for 2024-07-12
# ---------------------------------------------------------------------------
$cur = $wb.Worksheets.Item("current")

# The date label must stay literal text (not get auto-parsed into a date
# serial number) - a leading apostrophe forces Excel to store it as text.
$cur.Cells.Item(6, 1).Value = "'2024-07-12"
$cur.Cells.Item(6, 2).Value = 4
$cur.Cells.Item(6, 3).Value = 1
$cur.Cells.Item(6, 4).Value = 1
$cur.Cells.Item(6, 5).Value = 2

# ---------------------------------------------------------------------------
# 2) Add the new daily sheet "2024-07-12" after the last existing sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "2024-07-12"

# Header row
$newSheet.Cells.Item(1, 1).Value = "ID"
$newSheet.Cells.Item(1, 2).Value = "Время"
$newSheet.Cells.Item(1, 3).Value = "ФИО пациента"
$newSheet.Cells.Item(1, 4).Value = "Врач"
$newSheet.Cells.Item(1, 5).Value = "Врач_Индекс"
$newSheet.Cells.Item(1, 6).Value = "М\Ж\Р"
$newSheet.Cells.Item(1, 7).Value = "Дата рождения"
$newSheet.Cells.Item(1, 8).Value = "Причина"
$newSheet.Cells.Item(1, 9).Value = "Давление"

# Data rows - each is one filled-in visit form.
# ID(text) | Время(datetime) | ФИО(text) | Врач(text) | Врач_Индекс(number) |
# М\Ж\Р(text) | Дата рождения(text) | Причина(text) | Давление(text)
$rows = @(
  @{ id="1"; time=45485.24525712963; fio="iurhosthk";    doc="Karp_Kuzmin";     idx=3; sex="М"; dob="2006-07-04"; reason="рототщто";      pressure="7890" },
  @{ id="2"; time=45485.25785924769; fio="dxdjkughi";     doc="Karp_Kuzmin";     idx=3; sex="Р"; dob="2024-07-09"; reason="jgfghol";        pressure="4689" },
  @{ id="3"; time=45485.81736980324; fio="jykugk";        doc="Karp_Kuzmin";     idx=3; sex="Ж"; dob="2006-07-03"; reason="dstfui";         pressure="7890" },
  @{ id="1"; time=45485.83415211806; fio="patient name";  doc="Yefrem_Lebedev";  idx=2; sex="М"; dob="2006-07-03"; reason="reason beseda";  pressure="pressure" },
  @{ id="1"; time=45485.83810637737; fio="patient test";  doc="Desya_Osipov";    idx=5; sex="Ж"; dob="2006-07-11"; reason="beseda osipov";  pressure="pressure test" }
)

$r = 2
foreach ($row in $rows) {
  # A: ID - numeric-looking, force as text with a leading apostrophe
  $newSheet.Cells.Item($r, 1).Value = "'" + $row.id

  # B: Время - real date/time number with the custom datetime format
  $newSheet.Cells.Item($r, 2).Value = $row.time
  $newSheet.Cells.Item($r, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"

  # C: ФИО пациента (plain text)
  $newSheet.Cells.Item($r, 3).Value = $row.fio

  # D: Врач (plain text)
  $newSheet.Cells.Item($r, 4).Value = $row.doc

  # E: Врач_Индекс - an actual number
  $newSheet.Cells.Item($r, 5).Value = $row.idx

  # F: М\Ж\Р (plain text)
  $newSheet.Cells.Item($r, 6).Value = $row.sex

  # G: Дата рождения - date-looking text, force as text
  $newSheet.Cells.Item($r, 7).Value = "'" + $row.dob

  # H: Причина (plain text)
  $newSheet.Cells.Item($r, 8).Value = $row.reason

  # I: Давление - sometimes numeric-looking, force as text
  $newSheet.Cells.Item($r, 9).Value = "'" + $row.pressure

  $r = $r + 1
}

# Restore the originally active sheet ("2024-07-10" tab stays the active one).
$wb.Worksheets.Item("2024-07-10").Activate()
